# Add a new column AA ("SPY_Return") holding the running cumulative sum of
# column B ("SPY_ret"), mirroring the way column Z ("Total_PL") already
# accumulates the strategy's monthly P/L.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell AA1: copy the existing header formatting from Z1, then set the text.
$ws.Range("Z1").Copy()
$ws.Range("AA1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("AA1").Value = "SPY_Return"

# Data cells AA2:AA279: running total of B2:B279.
$ws.Range("AA2:AA279").Formula = "=SUM(`$B`$2:B2)"

# Convert the formulas to plain static values, matching the source data.
$ws.Range("AA2:AA279").Copy()
$ws.Range("AA2:AA279").PasteSpecial(-4163)   # xlPasteValues

$excel.CutCopyMode = 0
